$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.393.02"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.451.15"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.73"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.01"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.447.06"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  -0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.357"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.13"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.876.30"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.996.30"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.443.90"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.31"
$ws.Range("E20").Value = "  +5.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.36"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.11"
$ws.Range("E23").Value = "  +15.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.49"
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "619.41"
$ws.Range("E26").Value = "  +6.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  +5.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000104"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.561.19"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("E35").Value = "  +8.54%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.79"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +11.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "145.82"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.93"
$ws.Range("E45").Value = "  +0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "148.94"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.77"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.26"
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0536"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.602"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0234"
$ws.Range("E51").Value = "  +1.07%  "
